$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 41670064
$ws.Range("J70").Value = 41672044
$ws.Range("L70").Value = 125016132
$ws.Range("N70").Value = -125016672
$ws.Range("H73").Value = 41670064
$ws.Range("J73").Value = 41672044
$ws.Range("L73").Value = 125016132
$ws.Range("N73").Value = -125018004
$ws.Range("H132").Value = 2215.8572
$ws.Range("I132").Value = 2267.1853
$ws.Range("K132").Value = 6801.5559
$ws.Range("M132").Value = -4271.5559
$ws.Range("H137").Value = 3128.1892
$ws.Range("I137").Value = 3353.6667
$ws.Range("J137").Value = 2914.5789
$ws.Range("K137").Value = 10061.0001
$ws.Range("L137").Value = 8743.736699999999
$ws.Range("M137").Value = -7511.000100000001
$ws.Range("N137").Value = -13843.7367
$ws.Range("H138").Value = 1927791.1
$ws.Range("J138").Value = 2504763
$ws.Range("L138").Value = 7514289
$ws.Range("N138").Value = -7524569
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1764383.4
$ws.Range("I32").Value = 2086628.2
$ws.Range("J32").Value = 6684.5454
$ws.Range("K32").Value = 2086628.2
$ws.Range("L32").Value = 6684.5454
$ws.Range("M32").Value = -2086341.2
$ws.Range("N32").Value = -7258.5454
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H61").Value = 6364.6665
$ws.Range("I61").Value = 3467.5312
$ws.Range("K61").Value = 3467.5312
$ws.Range("M61").Value = -3255.5312
$ws.Range("H74").Value = 31958.143
$ws.Range("I74").Value = 38760.555
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 38760.555
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -37886.555
$ws.Range("N74").Value = -10748
$ws.Range("H77").Value = 31958.143
$ws.Range("I77").Value = 38760.555
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 193802.775
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -189434.775
$ws.Range("N77").Value = -53736
$ws.Range("H102").Value = 2317.3333
$ws.Range("I102").Value = 2317.3333
$ws.Range("K102").Value = 2317.3333
$ws.Range("M102").Value = -695.3332999999998
$ws.Range("H132").Value = 1823748.6
$ws.Range("I132").Value = 5273286.5
$ws.Range("J132").Value = 8202.263000000001
$ws.Range("K132").Value = 15819859.5
$ws.Range("L132").Value = 24606.789
$ws.Range("M132").Value = -15817329.5
$ws.Range("N132").Value = -29666.789
$ws.Range("H136").Value = 6364.6665
$ws.Range("I136").Value = 3467.5312
$ws.Range("K136").Value = 10402.5936
$ws.Range("M136").Value = -7852.5936
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3836.6667
$ws.Range("I105").Value = 2687.3333
$ws.Range("J105").Value = 5369.1113
$ws.Range("K105").Value = 2687.3333
$ws.Range("L105").Value = 5369.1113
$ws.Range("M105").Value = -940.3332999999998
$ws.Range("N105").Value = -8863.1113
$ws.Range("I107").Value = 56252070
$ws.Range("K107").Value = 56252070
$ws.Range("M107").Value = -56250150
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8713.462
$ws.Range("J31").Value = 10961.208
$ws.Range("L31").Value = 10961.208
$ws.Range("N31").Value = -11551.208
$ws.Range("H34").Value = 8713.462
$ws.Range("J34").Value = 10961.208
$ws.Range("L34").Value = 10961.208
$ws.Range("N34").Value = -11365.208
$ws.Range("H38").Value = 36446.5
$ws.Range("J38").Value = 36446.5
$ws.Range("L38").Value = 36446.5
$ws.Range("N38").Value = -37200.5
$ws.Range("H46").Value = 36446.5
$ws.Range("J46").Value = 36446.5
$ws.Range("L46").Value = 36446.5
$ws.Range("N46").Value = -36868.5
$ws.Range("H58").Value = 5864.55
$ws.Range("I58").Value = 1966.1428
$ws.Range("J58").Value = 10173.315
$ws.Range("K58").Value = 1966.1428
$ws.Range("L58").Value = 10173.315
$ws.Range("M58").Value = -1763.1428
$ws.Range("N58").Value = -10579.315
$ws.Range("H59").Value = 98181.67999999999
$ws.Range("J59").Value = 98181.67999999999
$ws.Range("L59").Value = 98181.67999999999
$ws.Range("N59").Value = -100471.68
$ws.Range("H68").Value = 38031.668
$ws.Range("J68").Value = 42047.5
$ws.Range("L68").Value = 42047.5
$ws.Range("N68").Value = -43545.5
$ws.Range("H71").Value = 38031.668
$ws.Range("J71").Value = 42047.5
$ws.Range("L71").Value = 126142.5
$ws.Range("N71").Value = -133630.5
$ws.Range("H74").Value = 125053620
$ws.Range("J74").Value = 61271
$ws.Range("L74").Value = 61271
$ws.Range("N74").Value = -63019
$ws.Range("H77").Value = 125053620
$ws.Range("J77").Value = 61271
$ws.Range("L77").Value = 183813
$ws.Range("N77").Value = -192549
$ws.Range("H132").Value = 5380.4546
$ws.Range("I132").Value = 2710.4707
$ws.Range("J132").Value = 8217.3125
$ws.Range("K132").Value = 8131.4121
$ws.Range("L132").Value = 24651.9375
$ws.Range("M132").Value = -5601.4121
$ws.Range("N132").Value = -29711.9375
$ws.Range("H134").Value = 5145.884
$ws.Range("I134").Value = 2817.2
$ws.Range("K134").Value = 8451.599999999999
$ws.Range("M134").Value = -5916.599999999999
$ws.Range("H136").Value = 5864.55
$ws.Range("I136").Value = 1966.1428
$ws.Range("J136").Value = 10173.315
$ws.Range("K136").Value = 5898.428400000001
$ws.Range("L136").Value = 30519.945
$ws.Range("M136").Value = -3348.428400000001
$ws.Range("N136").Value = -35619.945
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 956
$ws.Range("I64").Value = 412
$ws.Range("K64").Value = 1236
$ws.Range("M64").Value = -966
$ws.Range("H67").Value = 956
$ws.Range("I67").Value = 412
$ws.Range("K67").Value = 1236
$ws.Range("M67").Value = -300
$ws.Range("H68").Value = 2623.12
$ws.Range("I68").Value = 1999.875
$ws.Range("J68").Value = 2741.8333
$ws.Range("K68").Value = 5999.625
$ws.Range("L68").Value = 8225.499899999999
$ws.Range("M68").Value = -5188.625
$ws.Range("N68").Value = -9847.499899999999
$ws.Range("H71").Value = 2623.12
$ws.Range("I71").Value = 1999.875
$ws.Range("J71").Value = 2741.8333
$ws.Range("K71").Value = 17998.875
$ws.Range("L71").Value = 24676.4997
$ws.Range("M71").Value = -13942.875
$ws.Range("N71").Value = -32788.4997
$ws.Range("H74").Value = 2972
$ws.Range("J74").Value = 2972
$ws.Range("L74").Value = 8916
$ws.Range("N74").Value = -11038
$ws.Range("H77").Value = 2972
$ws.Range("J77").Value = 2972
$ws.Range("L77").Value = 26748
$ws.Range("N77").Value = -37356
$ws.Range("H102").Value = 10500
$ws.Range("J102").Value = 10500
$ws.Range("L102").Value = 31500
$ws.Range("N102").Value = -36368
$ws.Range("H131").Value = 3579.182
$ws.Range("J131").Value = 3637.1
$ws.Range("L131").Value = 10911.3
$ws.Range("N131").Value = -20991.3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 999
$ws.Range("I43").Value = 999
$ws.Range("K43").Value = 999
$ws.Range("M43").Value = -848
$ws.Range("H46").Value = 34499.5
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50312
$ws.Range("H126").Value = 6814.095
$ws.Range("I126").Value = 4799
$ws.Range("J126").Value = 7288.2354
$ws.Range("K126").Value = 14397
$ws.Range("L126").Value = 21864.7062
$ws.Range("M126").Value = -11927
$ws.Range("N126").Value = -26804.7062
$ws.Range("H132").Value = 8423.714
$ws.Range("I132").Value = 1988.3334
$ws.Range("J132").Value = 13250.25
$ws.Range("K132").Value = 5965.0002
$ws.Range("L132").Value = 39750.75
$ws.Range("M132").Value = -3435.0002
$ws.Range("N132").Value = -44810.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 15013
$ws.Range("J39").Value = 11688.333
$ws.Range("L39").Value = 11688.333
$ws.Range("N39").Value = -12608.333
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 26258844
$ws.Range("I81").Value = 2001649.6
$ws.Range("K81").Value = 4003299.2
$ws.Range("M81").Value = -4002238.2
$ws.Range("H84").Value = 26258844
$ws.Range("I84").Value = 2001649.6
$ws.Range("K84").Value = 20016496
$ws.Range("M84").Value = -20011192
$ws.Range("H136").Value = 24393486
$ws.Range("I136").Value = 41668360
$ws.Range("J136").Value = 5426.647
$ws.Range("K136").Value = 125005080
$ws.Range("L136").Value = 16279.941
$ws.Range("M136").Value = -125002530
$ws.Range("N136").Value = -21379.941
